$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Förändrad" (changed/updated) date column C was refreshed for every
# existing data row (2-400) from 2023-09-23 (45192) to 2023-10-03 (45202).
$ws.Range("C2:C400").Value = 45202

# Append the new record row (401) for case "A 45723-2023".
$ws.Range("A401").Value = "A 45723-2023"

$ws.Range("B401").Value = 45195
$ws.Range("B401").NumberFormat = "YYYY-MM-DD"

$ws.Range("C401").Value = 45202
$ws.Range("C401").NumberFormat = "YYYY-MM-DD"

$ws.Range("D401").Value = "VÄSTERNORRLANDS LÄN"
$ws.Range("E401").Value = "TIMRÅ"

$ws.Range("G401").Value = 0.9
$ws.Range("H401").Value = 0
$ws.Range("I401").Value = 0
$ws.Range("J401").Value = 0
$ws.Range("K401").Value = 0
$ws.Range("L401").Value = 0
$ws.Range("M401").Value = 0
$ws.Range("N401").Value = 0
$ws.Range("O401").Value = 0
$ws.Range("P401").Value = 0
$ws.Range("Q401").Value = 0

# R column keeps the wrap-text style even though this row has no species list.
$ws.Range("R401").Value = ""
$ws.Range("R401").WrapText = $true

# Row 400 picks up an explicit default row height once row 401 exists below it.
$ws.Rows.Item(400).RowHeight = 15
